# Remove "3350" from the fig 4 ("1-Day PEG Schematic") text box.
#
# Target changes (per commit "Remove 3350 from fig 4 schematic"):
#   "1-day of 15% PEG 3350 in drinking water (N = 6) "  -> "1-day of 15% PEG in drinking water (N = 6) "
#   "1-day PEG 3350 + 1-day recovery (N = 6)"           -> "1-day PEG + 1-day recovery (N = 6)"
#
# (The same diff also touches several datetimeFigureOut fields on other
# slides of the original multi-slide deck; this presentation only contains
# the single "Groups" schematic slide, which has no live date fields, so
# only the text edit below applies here.)

$p = $ppt.ActivePresentation

function Remove-SubstringFromShape($shp, [string]$needle) {
    if (-not $shp.HasTextFrame) { return }
    $tr = $shp.TextFrame.TextRange
    $text = $tr.Text
    if ($text -eq $null) { return }

    # Work from the right-most match backwards so earlier character
    # offsets are not invalidated by the deletion of later ones.
    $positions = New-Object System.Collections.ArrayList
    $searchFrom = 0
    while ($true) {
        $idx = $text.IndexOf($needle, $searchFrom)
        if ($idx -lt 0) { break }
        [void]$positions.Add($idx)
        $searchFrom = $idx + $needle.Length
    }

    for ($k = $positions.Count - 1; $k -ge 0; $k--) {
        $zeroBased = $positions[$k]
        $oneBased = $zeroBased + 1
        $tr.Characters($oneBased, $needle.Length).Text = ""
    }
}

function Process-Shapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        $isGroup = $false
        try { $isGroup = ($shp.Type -eq 6) } catch { $isGroup = $false }

        if ($isGroup) {
            Process-Shapes $shp.GroupItems
        } else {
            $hasText = $false
            try { $hasText = $shp.HasTextFrame } catch { $hasText = $false }
            if ($hasText) {
                $t = $shp.TextFrame.TextRange.Text
                if ($t -ne $null -and $t.Contains("3350")) {
                    Remove-SubstringFromShape $shp " 3350"
                }
            }
        }
    }
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    Process-Shapes $slide.Shapes
}
